$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.563.27"
$ws.Range("E2").Value = "'  +5.52%  "

$ws.Range("D3").Value = "'3.173.06"
$ws.Range("E3").Value = "'  +2.13%  "

$ws.Range("E4").Value = "'  +0.07%  "

$ws.Range("D5").Value = "'401.18"
$ws.Range("E5").Value = "'  +3.18%  "

$ws.Range("D6").Value = "'109.93"
$ws.Range("E6").Value = "'  +6.56%  "

$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "'  +1.27%  "

$ws.Range("E8").Value = "'  -0.05%  "

$ws.Range("E9").Value = "'  +4.16%  "

$ws.Range("D10").Value = "'39.15"
$ws.Range("E10").Value = "'  +4.81%  "

$ws.Range("D11").Value = "'0.0900"
$ws.Range("E11").Value = "'  +4.48%  "

$ws.Range("E12").Value = "'  +1.69%  "

$ws.Range("D13").Value = "'3.677.89"
$ws.Range("E13").Value = "'  +2.34%  "

$ws.Range("B14").Value = "'Chainlink"
$ws.Range("C14").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'19.08"
$ws.Range("E14").Value = "'  +2.14%  "

$ws.Range("B15").Value = "'Polkadot"
$ws.Range("C15").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'8.07"
$ws.Range("E15").Value = "'  +2.01%  "

$ws.Range("E16").Value = "'  +7.06%  "

$ws.Range("D17").Value = "'3.174.24"
$ws.Range("E17").Value = "'  +2.50%  "

$ws.Range("D18").Value = "'10.56"
$ws.Range("E18").Value = "'  -3.43%  "

$ws.Range("D19").Value = "'54.438.99"
$ws.Range("E19").Value = "'  +5.15%  "

$ws.Range("D20").Value = "'3.29"
$ws.Range("E20").Value = "'  +2.97%  "

$ws.Range("D22").Value = "'12.95"
$ws.Range("E22").Value = "'  +4.01%  "

$ws.Range("D23").Value = "'72.14"
$ws.Range("E23").Value = "'  +2.99%  "

$ws.Range("D24").Value = "'275.50"
$ws.Range("E24").Value = "'  +2.84%  "

$ws.Range("D25").Value = "'3.26"
$ws.Range("E25").Value = "'  +3.88%  "

$ws.Range("E26").Value = "'  -0.71%  "

$ws.Range("D27").Value = "'7.65"
$ws.Range("E27").Value = "'  +7.18%  "

$ws.Range("D28").Value = "'27.84"
$ws.Range("E28").Value = "'  +2.54%  "

$ws.Range("E29").Value = "'  -0.19%  "

$ws.Range("E31").Value = "'  +2.22%  "

$ws.Range("D32").Value = "'11.12"
$ws.Range("E32").Value = "'  +7.37%  "

$ws.Range("E33").Value = "'  +13.02%  "

$ws.Range("D34").Value = "'36.64"
$ws.Range("E34").Value = "'  +3.15%  "

$ws.Range("E35").Value = "'  +1.10%  "

$ws.Range("D36").Value = "'51.25"
$ws.Range("E36").Value = "'  +1.89%  "

$ws.Range("E38").Value = "'  -0.04%  "

$ws.Range("D39").Value = "'2.89"
$ws.Range("E39").Value = "'  +10.70%  "

$ws.Range("D40").Value = "'4.08"
$ws.Range("E40").Value = "'  +9.97%  "

$ws.Range("D41").Value = "'0.293"
$ws.Range("E41").Value = "'  +1.30%  "

$ws.Range("E42").Value = "'  +2.35%  "

$ws.Range("B43").Value = "'Celestia"
$ws.Range("C43").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'17.26"
$ws.Range("E43").Value = "'  +2.67%  "

$ws.Range("B44").Value = "'Monero"
$ws.Range("C44").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'132.27"
$ws.Range("E44").Value = "'  +2.57%  "

$ws.Range("E45").Value = "'  +1.19%  "

$ws.Range("E46").Value = "'  -0.20%  "

$ws.Range("D47").Value = "'2.48"
$ws.Range("E47").Value = "'  -0.61%  "

$ws.Range("D48").Value = "'2.08"
$ws.Range("E48").Value = "'  -0.48%  "

$ws.Range("D49").Value = "'2.102.66"
$ws.Range("E49").Value = "'  +2.82%  "

$ws.Range("D50").Value = "'0.0517"
$ws.Range("E50").Value = "'  +13.72%  "

$ws.Range("E51").Value = "'  +6.24%  "
